# Updating assignment 4 grades
# Adds a new "Assignment 4" column (G) with per-student grades, fixes a
# previously-miskeyed "Assignment 2" grade for one student, and moves the
# active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Assignment 4" column ------------------------------------

# Header cell: copy the formatting used by the other header cells (e.g. F1)
# and then set the text.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Assignment 4"

# "maximum grade" row: copy the formatting used by F2.
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("G2").Value = 100

# Data rows (3-43): copy the plain numeric formatting used in column A.
$ws.Range("A3").Copy()
$ws.Range("G3:G43").PasteSpecial(-4122)

$grades = @(50,0,100,100,100,75,100,100,100,100,0,100,0,100,75,100,100,100,100,93.75,100,100,100,100,100,100,100,100,100,100,0,100,100,0,100,100,100,100,100,100,0)
for ($i = 0; $i -lt $grades.Length; $i++) {
    $row = $i + 3
    $ws.Cells.Item($row, 7).Value = $grades[$i]
}

$excel.CutCopyMode = $false

# Give the new column a sensible display width (matches the others visually).
$ws.Columns.Item(7).ColumnWidth = 16.83

# --- Fix a miskeyed grade ---------------------------------------------------

$ws.Range("E22").Value = 88.24

# --- Update the saved selection --------------------------------------------

[void]$ws.Range("C5").Select()
